$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row,
# matching the refreshed figures from the source feed.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.126.30"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.827.65"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.67%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("E6").Value = "  -0.52%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4574"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +7.49%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3741"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.09%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07331"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.85%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8615"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.32%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.98"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.840.89"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.695"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.97%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "92.97"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +5.81%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.346"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("E16").Value = "  -0.22%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.63%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008833"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("E19").Value = "  -0.55%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "15.02"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "27.191.86"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.198"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "11.03"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.68%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.998"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.57%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "152.17"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.240"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +5.91%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.62"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.75%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "5.265"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.79%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "117.36"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.46%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.08864"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.04%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.194"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.21%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7582"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.971"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.76%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.469"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "

$ws.Range("E35").Value = "  -0.55%  "

$ws.Range("E36").Value = "  -0.71%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01969"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.99%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05289"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.98%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.5371"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.31%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "7.225"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.32%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.887"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1711"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.39%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.5227"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +11.60%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.626"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "10.74"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.34%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.964"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +9.31%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "106.38"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.11%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.677"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.36%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06358"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.9252"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.75%  "
